# "present the sales report data by state and to exclude all corporate information"
# -> Sort the Financial Data table by State (column B), add an AutoFilter over
#    the whole table (so the Corporate office type can be filtered out), and
#    leave the selection where the user ended up after doing the work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financial Data")
$ws.Select()

# Sort the data rows (A2:H20) by State (column B), ascending, keeping the
# header row out of the sort.
$sortRange = $ws.Range("A1:H20")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B20"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# The row reorder invalidates the previous shared-formula grouping; rewrite
# each Operating Profit / Net Profit formula explicitly so every cell holds
# its own (non-shared) formula, matching what Excel does after a sort.
for ($r = 2; $r -le 20; $r++) {
    $ws.Range("F$r").Formula = "=D$r-E$r"
    $ws.Range("H$r").Formula = "=F$r-G$r"
}

# Turn on AutoFilter across the full used range of the sheet so Corporate
# entries (and anything else) can be excluded interactively.
$fullRange = $ws.Range("A1:S27")
$fullRange.AutoFilter()

# AutoFilter defines a hidden sheet-scoped _FilterDatabase name that points at
# the filtered range.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Financial Data'!`$A`$1:`$S`$27")
$filterName.Visible = $false

# Leave the cursor where it ended up after the edit.
$ws.Range("I25").Select()
